$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlCenter = -4108

# --- Header row: B3 ("Points") gets centered style (s=2) ---
$ws.Range("B3").HorizontalAlignment = $xlCenter

# --- Row 4: Andriy Zhidko ---
$ws.Range("B4").Value = "?"
$ws.Range("B4").HorizontalAlignment = $xlCenter
$ws.Range("D4").Value = "+"
$ws.Range("D4").HorizontalAlignment = $xlCenter

# --- Row 5: Mark Vlasenko ---
$ws.Range("B5").Value = 0
$ws.Range("B5").HorizontalAlignment = $xlCenter
$ws.Range("D5").Value = "-"
$ws.Range("D5").HorizontalAlignment = $xlCenter

# --- Row 6: Ivan Maidaniuk ---
$ws.Range("B6").Value = "?"
$ws.Range("B6").HorizontalAlignment = $xlCenter
$ws.Range("D6").Value = "+"
$ws.Range("D6").HorizontalAlignment = $xlCenter

# --- Row 7: Volodymyr Tarasyuk (name becomes bold) ---
$ws.Range("A7").Font.Bold = $true
$ws.Range("B7").Value = 8
$ws.Range("B7").HorizontalAlignment = $xlCenter
$ws.Range("D7").Value = "+"
$ws.Range("D7").HorizontalAlignment = $xlCenter

# --- Row 8: Nataliia Striukova ---
$ws.Range("B8").Value = 8
$ws.Range("B8").HorizontalAlignment = $xlCenter
$ws.Range("D8").Value = "+"
$ws.Range("D8").HorizontalAlignment = $xlCenter

# --- Row 9: Anatolii Donchenko ---
$ws.Range("B9").Value = 4
$ws.Range("B9").HorizontalAlignment = $xlCenter
$ws.Range("D9").Value = "+"
$ws.Range("D9").HorizontalAlignment = $xlCenter

# --- Row 10: Veronika Pliuslina ---
$ws.Range("B10").Value = 4
$ws.Range("B10").HorizontalAlignment = $xlCenter
$ws.Range("D10").Value = "+"
$ws.Range("D10").HorizontalAlignment = $xlCenter

# --- Row 11: Volodymyr Lezhnenko ---
$ws.Range("B11").Value = 4
$ws.Range("B11").HorizontalAlignment = $xlCenter
$ws.Range("D11").Value = "+"
$ws.Range("D11").HorizontalAlignment = $xlCenter

# --- Row 12: Oleg Tokar ---
$ws.Range("B12").Value = 4
$ws.Range("B12").HorizontalAlignment = $xlCenter
$ws.Range("D12").Value = "+"
$ws.Range("D12").HorizontalAlignment = $xlCenter

# --- Row 13: Olga Kovalyshyn ---
$ws.Range("B13").Value = 4
$ws.Range("B13").HorizontalAlignment = $xlCenter
$ws.Range("D13").Value = "+"
$ws.Range("D13").HorizontalAlignment = $xlCenter

# --- Row 14: Ihor Horovetskyi ---
$ws.Range("B14").Value = 4
$ws.Range("B14").HorizontalAlignment = $xlCenter
$ws.Range("D14").Value = "+"
$ws.Range("D14").HorizontalAlignment = $xlCenter

# --- Row 15: Khrystyna Lysiuk (B15 stays empty but styled) ---
$ws.Range("B15").HorizontalAlignment = $xlCenter
$ws.Range("D15").Value = "-"
$ws.Range("D15").HorizontalAlignment = $xlCenter

# --- Row 16: Andrii Vishchansky ---
$ws.Range("B16").Value = 4
$ws.Range("B16").HorizontalAlignment = $xlCenter
$ws.Range("D16").Value = "+"
$ws.Range("D16").HorizontalAlignment = $xlCenter

# --- Row 17: brand-new blank spacer row, just an empty styled B cell ---
$ws.Range("B17").HorizontalAlignment = $xlCenter

# --- Row 18: "Friday Group" header ---
$ws.Range("B18").HorizontalAlignment = $xlCenter

# --- Row 19: Alina Keda ---
$ws.Range("B19").HorizontalAlignment = $xlCenter

# --- Row 20: Tetyana Suprunova ---
$ws.Range("B20").HorizontalAlignment = $xlCenter

# --- Row 21: Yaroslav Danylchenko - Points becomes a formula ---
$ws.Range("B21").Formula = "=5+8"

# --- Row 22: Evgen Bruhov ---
$ws.Range("B22").HorizontalAlignment = $xlCenter

# --- Row 23: Olena Sudarkina ---
$ws.Range("B23").HorizontalAlignment = $xlCenter

# --- Row 24: Maxim Chireychik ---
$ws.Range("B24").Value = 8
$ws.Range("B24").HorizontalAlignment = $xlCenter

# --- Row 25: Mykhailo Sapiegin ---
$ws.Range("B25").HorizontalAlignment = $xlCenter

# --- Row 26: Andrii Rasskazov ---
$ws.Range("B26").HorizontalAlignment = $xlCenter

# --- Row 27: Mykhailo Semenikhin ---
$ws.Range("B27").Value = "?"
$ws.Range("B27").HorizontalAlignment = $xlCenter

# --- Row 28: Stanislav Maryenko ---
$ws.Range("B28").HorizontalAlignment = $xlCenter

# --- Row 29: Oleksandr Fedan ---
$ws.Range("B29").HorizontalAlignment = $xlCenter

# --- Row 30: Andriy Bilotskyy ---
$ws.Range("B30").HorizontalAlignment = $xlCenter

# --- Selection moved to A34, right below the last data row, ahead of checking the Friday group ---
$ws.Range("A34").Select()
